# Apply "Updated symbol list on Wed Dec 28 11:54:03 UTC 2022 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text (non-numeric-looking) cell updates: Coin name / Link / composite label columns ---
$ws.Range("E18").Value = "17OneONEWorstin24h"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("E41").Value = "40KickTokenKICK"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("E42").Value = "41BKEXTokenBKK"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("E43").Value = "42CEJICEJI"

# --- Price column (D) updates. These cells hold numeric-looking text that must stay as text,
# so force the "@" text format before assigning, then drop back to the default (unstyled) cell
# style afterwards so no stray number format sticks around on the cell.
$priceUpdates = @(
    , @("D3", "23.88")
    , @("D4", "5.256")
    , @("D5", "0.05822")
    , @("D6", "6.461")
    , @("D7", "3.330")
    , @("D8", "0.8084")
    , @("D9", "0.8733")
    , @("D10", "0.1381")
    , @("D11", "0.07267")
    , @("D12", "0.03084")
    , @("D13", "0.03054")
    , @("D14", "0.09330")
    , @("D15", "3.850")
    , @("D16", "0.001555")
    , @("D17", "0.04698")
    , @("D18", "0.0006017")
    , @("D19", "0.006141")
    , @("D20", "0.001262")
    , @("D21", "0.004595")
    , @("D22", "0.00008698")
    , @("D23", "3.561")
    , @("D24", "2.178")
    , @("D40", "0.03775")
    , @("D41", "0.006288")
    , @("D42", "0.1053")
    , @("D43", "0.002399")
    , @("D44", "0.007964")
    , @("D45", "0.00005522")
    , @("D47", "0.5976")
    , @("D48", "0.01416")
)
foreach ($pair in $priceUpdates) {
    $cell = $ws.Range($pair[0])
    $cell.NumberFormat = "@"
    $cell.Value = $pair[1]
    $cell.Style = "Normal"
}
